$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy header style (bold, border, centered) from an existing header cell (e.g. F1) to G1:H1
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Update existing row 2 values
$ws.Range("B2").Value = 0.1092714316659435
$ws.Range("C2").Value = 0.9980038697206826
$ws.Range("D2").Value = 0.2430551190367126
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor(n_estimators=100))])"

# New row 2 values
$ws.Range("G2").Value = 0.1311458841167526
$ws.Range("H2").Value = 0.991
